# Rename testdateN -> test_dateN for N = 1..9 in the equipment template.
#
# The target OOXML keeps each occurrence split across runs:
#   - the original leading "t" run is left untouched
#   - the original "estdate" run is shortened to a prefix ("est" for most
#     occurrences, "es" for the 2nd occurrence) and two brand-new runs carry
#     the remaining text (either "_"+"date" or "t_"+"date")
#   - the original trailing number run is left untouched
#
# Word's text-replace APIs rebuild/merge the whole paragraph's runs into a
# single run whenever the paragraph's text content changes, so we first do
# the textual rename (accepting the merge into one run per paragraph) and
# then re-split that merged run into the required pieces using a harmless
# Bold-on/Bold-off toggle, which this engine treats as a pure formatting
# boundary that does not re-merge the paragraph.

$d = $word.ActiveDocument

# Step 1: perform the textual rename for every occurrence. Each paragraph's
# runs collapse into a single run containing the full "test_dateN" text.
for ($n = 1; $n -le 9; $n++) {
    $old = "testdate" + $n
    $new = "test_date" + $n
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Step 2: re-split each merged "test_dateN" run back into the pieces that
# match the target document structure.
for ($n = 1; $n -le 9; $n++) {
    $needle = "test_date" + $n
    $rng = $d.Content
    $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $s = $rng.Start

    if ($n -eq 2) {
        # "t" | "es" | "t_" | "date" | "2"
        $cut1 = 1
        $cut2 = 3
        $cut3 = 5
        $cut4 = 9
    } else {
        # "t" | "est" | "_" | "date" | "N"
        $cut1 = 1
        $cut2 = 4
        $cut3 = 5
        $cut4 = 9
    }

    $p1 = $d.Range($s, $s + $cut1)
    $p1.Bold = 1
    $p1.Bold = 0

    $p2 = $d.Range($s + $cut1, $s + $cut2)
    $p2.Bold = 1
    $p2.Bold = 0

    $p3 = $d.Range($s + $cut2, $s + $cut3)
    $p3.Bold = 1
    $p3.Bold = 0

    $p4 = $d.Range($s + $cut3, $s + $cut4)
    $p4.Bold = 1
    $p4.Bold = 0
}
